$d = $word.ActiveDocument

function InsertRun($text, $super) {
    $startPos = $d.Content.End
    $r = $d.Range($startPos, $startPos)
    $r.InsertAfter($text)
    $r.Font.Name = "SolaimanLipi"
    $r.Font.NameBi = "SolaimanLipi"
    if ($super) {
        $r.Font.Superscript = $true
    }
}

function InsertPara() {
    $startPos = $d.Content.End
    $r = $d.Range($startPos, $startPos)
    $r.InsertParagraphAfter()
}

# Step 1: extend the final existing run text (delete + fresh insert to avoid run-merge with preceding "-" run)
$searchRange = $d.Content.Duplicate
$searchRange.Find.Execute("আকৃতির মৌলিক সীমা আছে। ", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$targetStart = $searchRange.Start
$targetEnd = $searchRange.End
$delRange = $d.Range($targetStart, $targetEnd)
$delRange.Delete()
$insRange = $d.Range($targetStart, $targetStart)
$insRange.InsertAfter("আকৃতির মৌলিক সীমা আছে। এ সীমাই স্থানের ক্ষেত্রফল ও আয়তনের কোয়ান্টা ঠিক করে দেয়। যা পরিমাপ করা হয় প্ল্যাঙ্ক দৈর্ঘ্যের মাধ্যমে। প্ল্যাঙ্ক দৈর্ঘ্য ১")
$insRange.Font.Name = "SolaimanLipi"
$insRange.Font.NameBi = "SolaimanLipi"

# Step 2: remaining runs of paragraph (part 0)
InsertRun "." $false
InsertRun "৬ " $false
InsertRun "× " $false
InsertRun "১০" $false
InsertRun "-" $false
InsertRun "৩৩" $true
InsertRun " মিটার। যা প্রোটনের ব্যাসের দশ লক্ষ" $false
InsertRun "-" $false
InsertRun "কোটি" $false
InsertRun "-" $false
InsertRun "কোটি ভাগের এক ভাগের সমান। " $false

# Step 3: empty paragraph
InsertPara

# Step 4: new paragraph with content (part 2)
InsertPara
InsertRun "ভিন্ন ভিন্ন স্পিন নেটওয়ার্ক লুপগুলোকে ভিন্ন ভিন্নভাবে জোড়া দেয়। স্থানের আকার" $false
InsertRun "-" $false
InsertRun "আকৃতির ভিন্ন ভিন্ন কোয়ান্টা অবস্থাও এভাবে তৈরি হয়। স্পিন নেটওয়ার্কের বিবর্তন " $false
InsertRun "(" $false
InsertRun "এক আকৃতির সঙ্গে অন্য আকৃতির পরিবর্তনশীল সম্পর্ক" $false
InsertRun ") " $false
InsertRun "থেকে জন্ম হয় স্পিনফোমের। সুপারপজিশন নামে একটি জিনিসের মধ্যে স্পিনফোমের সংযোজোনের মাধ্যমে উদীয়মান স্থান" $false
InsertRun "-" $false
InsertRun "কালের ব্যাখ্যা পাওয়া যায়। স্থান" $false
InsertRun "-" $false
InsertRun "কালের এ কাঠামো কোয়ান্টাম পদার্থবিদ্যার সাথে সহ" $false
InsertRun "-" $false
InsertRun "উৎপন্ন হয়। " $false
InsertRun "(" $false
InsertRun "পরিমাপ করার আগ পর্যন্ত একই সময়ে একটি কোয়ান্টাম সিস্টেম বহু অবস্থায় থাকতে পারে। এরই নাম সুপারপজিশন।" $false
InsertRun ") " $false

# Step 5: empty paragraph
InsertPara

# Step 6: new paragraph with content (part 4)
InsertPara
InsertRun "সংক্ষেপে এটাই লুপ কোয়ান্টাম গ্র্যাভিটি বা এলকিউজি। বর্তমানে " $false
InsertRun "(" $false
InsertRun "২০১৮ সালে এ বই লেখার সময়" $false
InsertRun ") " $false
InsertRun "এর বয়স ৩০ বছর। বর্তমানে সারা বিশ্বের ত্রিশটি গবেষণা দল আগ্রহের বস্তু এটি। আপেক্ষিকতা তত্ত্ব থেকে এখানে আসা সহজ ছিল না। পাড়ি দিতে হয়েছে চড়াই" $false
InsertRun "-" $false
InsertRun "উৎরাই। সামনে অনেক বাধা আছে এখনও। তার ওপর তত্ত্বটির গ্রহণযোগ্যতা পরীক্ষা করার উপায় বের করতে হবে। " $false
InsertRun "(" $false
InsertRun "পরীক্ষাযোগ্য না হলে কোনো তত্ত্বই বিজ্ঞানের অংশ হয়ে ওঠে না। থেকে যায় দর্শন।" $false
InsertRun ") " $false
